$d = $word.ActiveDocument

$frontHeading = $null
$firstParaOld = $null
foreach ($p in $d.Paragraphs) {
    if ($frontHeading -eq $null -and $p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.Trim() -eq "Frontmatter") {
        $frontHeading = $p
    } elseif ($frontHeading -ne $null -and $firstParaOld -eq $null -and $p.Style.NameLocal -eq "First Paragraph") {
        $firstParaOld = $p
        break
    }
}

# Insert a new blank paragraph right after the "Frontmatter" heading (before the FirstParagraph text)
$insertPoint = $d.Range($frontHeading.Range.End, $frontHeading.Range.End)
$insertPoint.InsertParagraphAfter()

# The new blank paragraph is now between the heading and the old text paragraph; find it
$abstractHeadingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $frontHeading.Range.End) {
        $abstractHeadingPara = $p
        break
    }
}
$abstractHeadingPara.Range.Text = "Abstract"
$abstractHeadingPara.Style = "Heading 3"

# Re-locate the old FirstParagraph text paragraph (positions shifted after insertion above)
$firstParaOld = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $abstractHeadingPara.Range.End) {
        $firstParaOld = $p
        break
    }
}
Write-Host "firstParaOld: '" $firstParaOld.Range.Text "'"

$abstractBodyStart = $firstParaOld.Range.Start

$fullText = 'Technologies including PCs, smartphones, and cloud computing have transformed the world: In our daily lives we interact with many businesses and public services who (in pursuit of cost-saving) increasingly seek to rely on data collection and processing rather than face-to-face user interactions to inform their decisions. This creates an imbalance of power between those who hold data and the individuals about whom data is stored, who cannot easily see their personal data or how it is used. This Digital Civics PhD research explores, from a pragmatic, constructivist perspective, the topic of Human Data Relations. Through two qualitative case studies across public and private sectors, it answers the question, “What do we want from data, and from those who hold data about us?”. Case Study One focuses on Early Help social care: Through four workshops with supported families, social workers and staff, a deep understanding of the individual perspective on civic personal data use is established. Shared data interaction is explored as a means to shift the balance of power towards the individual while maintaining an effective care relationship. Case Study Two is a three-month study exploring 10 participants’ experience of using GDPR data access rights to view their own data, resulting in insights into individual needs and the challenges of data-centric service relationships, and recommendations for improvement of policies and practices. With reference to literature from the fields of Personal Information Management, Human Data Interaction and MyData personal data ecosystems, these case studies contribute to a unified understanding of six core needs that people have in Human Data Relations. In the final chapter, the thesis discusses the practical pursuit of these goals, drawing on first-hand knowledge acquired from expert participation in industrial research projects at BBC R&D and Hestia.ai/SITRA, providing a workable roadmap for future research and innovation.'
$firstParaOld.Range.Text = $fullText

# Apply italic formatting to the designated sub-ranges (offsets relative to paragraph start)
$d.Range($abstractBodyStart + 339, $abstractBodyStart + 357).Font.Italic = 1
$d.Range($abstractBodyStart + 499, $abstractBodyStart + 513).Font.Italic = 1
$d.Range($abstractBodyStart + 596, $abstractBodyStart + 616).Font.Italic = 1
$d.Range($abstractBodyStart + 715, $abstractBodyStart + 782).Font.Italic = 1
$d.Range($abstractBodyStart + 810, $abstractBodyStart + 820).Font.Italic = 1
$d.Range($abstractBodyStart + 1002, $abstractBodyStart + 1025).Font.Italic = 1
$d.Range($abstractBodyStart + 1237, $abstractBodyStart + 1260).Font.Italic = 1
$d.Range($abstractBodyStart + 1498, $abstractBodyStart + 1529).Font.Italic = 1
$d.Range($abstractBodyStart + 1531, $abstractBodyStart + 1553).Font.Italic = 1
$d.Range($abstractBodyStart + 1558, $abstractBodyStart + 1564).Font.Italic = 1
$d.Range($abstractBodyStart + 1651, $abstractBodyStart + 1665).Font.Italic = 1
$d.Range($abstractBodyStart + 1755, $abstractBodyStart + 1772).Font.Italic = 1
$d.Range($abstractBodyStart + 1941, $abstractBodyStart + 1983).Font.Italic = 1

Write-Host "DONE_TEXT_AND_ITALICS"

# Add the "abstract" bookmark spanning the Abstract heading paragraph through the end of its body paragraph
$abstractBookmarkRange = $d.Range($abstractHeadingPara.Range.Start, $firstParaOld.Range.End)
$d.Bookmarks.Add("abstract", $abstractBookmarkRange)
Write-Host "DONE_BOOKMARK"
